# Update column G ("K") values on the active sheet to reflect the
# regenerated strike-count data (K instead of Strike#).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 2
    9  = 2
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 3
    15 = 1
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 2
    22 = 0
    23 = 1
    24 = 2
    25 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
